# Mise à jour de l'application
# Adds a new attendance column (AF) for 2025-08-21 (serial 45890) and
# corrects two retroactive attendance entries on row 13 (AD13/AE13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date header in AF1, formatted like the other date cells (AE1) ---
$ws.Range("AF1").Value = 45890
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)   # xlPasteFormats

# --- New day's attendance value per player row, formatted like AE<row> ---
$attendance = @{
    2  = "P"
    3  = "M"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    12 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "RH"
    17 = "B"
    18 = "P"
    19 = "RH"
    20 = "P"
    21 = "P"
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "RH"
}

foreach ($row in 2..27) {
    $value = $attendance[$row]
    $ws.Range("AF$row").Value = $value
    $ws.Range("AE$row").Copy()
    $ws.Range("AF$row").PasteSpecial(-4122)   # xlPasteFormats
}

# --- Retroactive correction on row 13: A -> B for AD13 and AE13 ---
$ws.Range("AD13").Value = "B"
$ws.Range("AE13").Value = "B"

# --- View state: move selection to match the new active cell ---
$ws.Range("AH24").Select()
